$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (village), shifting collector_name etc. right.
$ws.Range("E1").EntireColumn.Insert()

# Match the column width of the column to the left (admin_level_2), same as
# Excel's native "copy formatting from left" behaviour on column insert.
$ws.Range("E1").ColumnWidth = $ws.Range("D1").ColumnWidth

# Set the header for the newly inserted column.
$ws.Range("E1").Value = "village"
